# Rebuild the Ishan Kishan (Mumbai Indians) match log:
#  - row 2 becomes a new "Oct 31 2020 vs Delhi Capitals" entry
#  - the former row 2 shifts down to row 3 (unchanged values)
#  - rows 4-13 are additional historical matches
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns G:K (runs, balls, 4s, 6s, SR) are stored as text in the source data,
# so force a text number format before writing numeric-looking strings into them
# -- otherwise Excel auto-converts "72" etc. into real numbers.
$ws.Range("G2:K13").NumberFormat = "@"

$data = @(
    @(' Oct 31 2020', ' Dubai (DSC)', 'Mumbai won by 9 wickets (with 34 balls remaining)', 'Mumbai Indians', 'Delhi Capitals', 'Ishan Kishan ', '72', '47', '8', '3', '153.19'),
    @(' Nov 10 2020', ' Dubai (DSC)', 'Mumbai won by 5 wickets (with 8 balls remaining)', 'Mumbai Indians', 'Delhi Capitals', 'Ishan Kishan ', '33', '19', '3', '1', '173.68'),
    @(' Oct 6 2020', ' Abu Dhabi', 'Mumbai won by 57 runs', 'Mumbai Indians', 'Rajasthan Royals', 'Ishan Kishan ', '0', '1', '0', '0', '0.00'),
    @(' Oct 23 2020', ' Sharjah', 'Mumbai won by 10 wickets (with 46 balls remaining)', 'Mumbai Indians', 'Chennai Super Kings', 'Ishan Kishan ', '68', '37', '6', '5', '183.78'),
    @(' Oct 18 2020', ' Dubai (DSC)', 'Match tied (Kings XI won the one-over eliminator)', 'Mumbai Indians', 'Kings XI Punjab', 'Ishan Kishan ', '7', '7', '1', '0', '100.00'),
    @(' Oct 25 2020', ' Abu Dhabi', 'Royals won by 8 wickets (with 10 balls remaining)', 'Mumbai Indians', 'Rajasthan Royals', 'Ishan Kishan ', '37', '36', '4', '1', '102.77'),
    @(' Oct 4 2020', ' Sharjah', 'Mumbai won by 34 runs', 'Mumbai Indians', 'Sunrisers Hyderabad', 'Ishan Kishan ', '31', '23', '1', '2', '134.78'),
    @(' Oct 11 2020', ' Abu Dhabi', 'Mumbai won by 5 wickets (with 2 balls remaining)', 'Mumbai Indians', 'Delhi Capitals', 'Ishan Kishan ', '28', '15', '2', '2', '186.66'),
    @(' Nov 5 2020', ' Dubai (DSC)', 'Mumbai won by 57 runs', 'Mumbai Indians', 'Delhi Capitals', 'Ishan Kishan ', '55', '30', '4', '3', '183.33'),
    @(' Nov 3 2020', ' Sharjah', 'Sunrisers won by 10 wickets (with 17 balls remaining)', 'Mumbai Indians', 'Sunrisers Hyderabad', 'Ishan Kishan ', '33', '30', '1', '2', '110.00'),
    @(' Sep 28 2020', ' Dubai (DSC)', 'Match tied (RCB won the one-over eliminator)', 'Mumbai Indians', 'Royal Challengers Bangalore', 'Ishan Kishan ', '99', '58', '2', '9', '170.68'),
    @(' Oct 1 2020', ' Abu Dhabi', 'Mumbai won by 48 runs', 'Mumbai Indians', 'Kings XI Punjab', 'Ishan Kishan ', '28', '32', '1', '1', '87.50')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

